# Collapse the split runs in the title, author and abstract paragraphs
# into single runs holding the full (unchanged) text, mirroring a
# "refresh" pass that re-saved the document with merged runs.

$d = $word.ActiveDocument

# Title: "Factsheet:" " " "List" " " "of" " " "derivatives"
#   -> "Factsheet: List of derivatives"
$d.Content.Find.Execute(
    "Factsheet: List of derivatives", $true, $false, $false, $false, $false,
    $true, 1, $false, "Factsheet: List of derivatives", 2)

# Author: "Tom" " " "Coleman" -> "Tom Coleman"
$d.Content.Find.Execute(
    "Tom Coleman", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tom Coleman", 2)

# Abstract: "A" " " "list" " " "of" " " "common" " " "(and" " " "some" " "
#   "uncommon)" " " "derivatives" " " "of" " " "functions."
#   -> "A list of common (and some uncommon) derivatives of functions."
$d.Content.Find.Execute(
    "A list of common (and some uncommon) derivatives of functions.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A list of common (and some uncommon) derivatives of functions.", 2)
